# "Code integration and output" — rewire the fbsweep sheet down to a single
# Node column and add a new "output" sheet (Node / PU / Polar Form / Cartesian Form).

$wb = $excel.ActiveWorkbook

# --- fbsweep: drop the old Sweep 1..5 header columns, keep only "Node" ---
$fbsweep = $wb.Worksheets.Item("fbsweep")
$fbsweep.Range("A1").Value = "Node"
$fbsweep.Range("B1:F1").ClearContents()

# --- add the new "output" sheet after fbsweep (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$output = $wb.Worksheets.Add($null, $lastSheet)
$output.Name = "output"

$output.Range("A1").Value = "Node"
$output.Range("B1").Value = "PU"
$output.Range("C1").Value = "Polar Form"
$output.Range("D1").Value = "Cartesian Form"

$output.Columns.Item(3).ColumnWidth = 9.7
$output.Columns.Item(4).ColumnWidth = 13.6

# --- restore per-sheet cursor positions, finishing on "voltages" as active tab ---
$edges = $wb.Worksheets.Item("edges")
[void]$edges.Activate()
[void]$edges.Range("B29").Select()

$nodes = $wb.Worksheets.Item("nodes")
[void]$nodes.Activate()
[void]$nodes.Range("C20").Select()

$settings = $wb.Worksheets.Item("settings")
[void]$settings.Activate()
[void]$settings.Range("F30").Select()

[void]$fbsweep.Activate()
[void]$fbsweep.Range("B4").Select()

[void]$output.Activate()
[void]$output.Range("C6").Select()

$voltages = $wb.Worksheets.Item("voltages")
[void]$voltages.Activate()
[void]$voltages.Range("B1").Select()
